$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = 201
$ws.Cells.Item(2, 2).Value = 9
$ws.Cells.Item(2, 3).Value = 30
$ws.Cells.Item(2, 4).Value = 15
$ws.Cells.Item(2, 5).Value = 45
$ws.Cells.Item(2, 6).Value = 30

$ws.Cells.Item(3, 1).Value = 801
$ws.Cells.Item(3, 2).Value = 3
$ws.Cells.Item(3, 3).Value = 67
$ws.Cells.Item(3, 4).Value = 65
$ws.Cells.Item(3, 5).Value = 52
$ws.Cells.Item(3, 6).Value = 45

$ws.Cells.Item(4, 1).Value = 1203
$ws.Cells.Item(4, 2).Value = 3
$ws.Cells.Item(4, 3).Value = 15
$ws.Cells.Item(4, 4).Value = 15
$ws.Cells.Item(4, 5).Value = 15
$ws.Cells.Item(4, 6).Value = 15

$ws.Cells.Item(5, 1).Value = 101
$ws.Cells.Item(5, 2).Value = 9
$ws.Cells.Item(5, 3).Value = 30
$ws.Cells.Item(5, 4).Value = 15
$ws.Cells.Item(5, 5).Value = 60
$ws.Cells.Item(5, 6).Value = 15

$ws.Cells.Item(6, 1).Value = 401
$ws.Cells.Item(6, 2).Value = 9
$ws.Cells.Item(6, 3).Value = 48
$ws.Cells.Item(6, 4).Value = 67
$ws.Cells.Item(6, 5).Value = 75
$ws.Cells.Item(6, 6).Value = 45

$ws.Cells.Item(7, 1).Value = 1202
$ws.Cells.Item(7, 2).Value = 2
$ws.Cells.Item(7, 3).Value = 10
$ws.Cells.Item(7, 4).Value = 10
$ws.Cells.Item(7, 5).Value = 10
$ws.Cells.Item(7, 6).Value = 10

$ws.Cells.Item(8, 1).Value = 902
$ws.Cells.Item(8, 2).Value = 1
$ws.Cells.Item(8, 3).Value = 0
$ws.Cells.Item(8, 4).Value = 0
$ws.Cells.Item(8, 5).Value = 0
$ws.Cells.Item(8, 6).Value = 0

$ws.Cells.Item(9, 1).Value = 1001
$ws.Cells.Item(9, 2).Value = 18
$ws.Cells.Item(9, 3).Value = 30
$ws.Cells.Item(9, 4).Value = 75
$ws.Cells.Item(9, 5).Value = 60
$ws.Cells.Item(9, 6).Value = 72

$ws.Cells.Item(10, 1).Value = 301
$ws.Cells.Item(10, 2).Value = 6
$ws.Cells.Item(10, 3).Value = 45
$ws.Cells.Item(10, 4).Value = 30
$ws.Cells.Item(10, 5).Value = 60
$ws.Cells.Item(10, 6).Value = 45

$ws.Cells.Item(11, 1).Value = 501
$ws.Cells.Item(11, 2).Value = 9
$ws.Cells.Item(11, 3).Value = 52
$ws.Cells.Item(11, 4).Value = 30
$ws.Cells.Item(11, 5).Value = 75
$ws.Cells.Item(11, 6).Value = 45

$ws.Cells.Item(12, 1).Value = 601
$ws.Cells.Item(12, 2).Value = 9
$ws.Cells.Item(12, 3).Value = 60
$ws.Cells.Item(12, 4).Value = 67
$ws.Cells.Item(12, 5).Value = 60
$ws.Cells.Item(12, 6).Value = 42

$ws.Cells.Item(13, 1).Value = 1201
$ws.Cells.Item(13, 2).Value = 2
$ws.Cells.Item(13, 3).Value = 10
$ws.Cells.Item(13, 4).Value = 10
$ws.Cells.Item(13, 5).Value = 10
$ws.Cells.Item(13, 6).Value = 10

$ws.Cells.Item(14, 1).Value = 901
$ws.Cells.Item(14, 2).Value = 16
$ws.Cells.Item(14, 3).Value = 15
$ws.Cells.Item(14, 4).Value = 45
$ws.Cells.Item(14, 5).Value = 60
$ws.Cells.Item(14, 6).Value = 60

$ws.Cells.Item(15, 1).Value = 701
$ws.Cells.Item(15, 2).Value = 3
$ws.Cells.Item(15, 3).Value = 90
$ws.Cells.Item(15, 4).Value = 45
$ws.Cells.Item(15, 5).Value = 97
$ws.Cells.Item(15, 6).Value = 15

$ws.Cells.Item(16, 1).Value = 502
$ws.Cells.Item(16, 2).Value = 0
$ws.Cells.Item(16, 3).Value = 4
$ws.Cells.Item(16, 4).Value = 0
$ws.Cells.Item(16, 5).Value = 0
$ws.Cells.Item(16, 6).Value = 0

$ws.Cells.Item(17, 1).Value = 802
$ws.Cells.Item(17, 2).Value = 0
$ws.Cells.Item(17, 3).Value = 4
$ws.Cells.Item(17, 4).Value = 5
$ws.Cells.Item(17, 5).Value = 4
$ws.Cells.Item(17, 6).Value = 0

$ws.Cells.Item(18, 1).Value = 2
$ws.Cells.Item(18, 2).Value = 0
$ws.Cells.Item(18, 3).Value = 2
$ws.Cells.Item(18, 4).Value = 2
$ws.Cells.Item(18, 5).Value = 2
$ws.Cells.Item(18, 6).Value = 2

$ws.Cells.Item(19, 1).Value = 3
$ws.Cells.Item(19, 2).Value = 0
$ws.Cells.Item(19, 3).Value = 3
$ws.Cells.Item(19, 4).Value = 3
$ws.Cells.Item(19, 5).Value = 3
$ws.Cells.Item(19, 6).Value = 3

$ws.Cells.Item(20, 1).Value = 1101
$ws.Cells.Item(20, 2).Value = 0
$ws.Cells.Item(20, 3).Value = 15
$ws.Cells.Item(20, 4).Value = 30
$ws.Cells.Item(20, 5).Value = 30
$ws.Cells.Item(20, 6).Value = 0

$ws.Cells.Item(21, 1).Value = 1
$ws.Cells.Item(21, 2).Value = 0
$ws.Cells.Item(21, 3).Value = 2
$ws.Cells.Item(21, 4).Value = 2
$ws.Cells.Item(21, 5).Value = 2
$ws.Cells.Item(21, 6).Value = 2

